# Update cryptos price/volume table (GitHub Actions scheduled refresh).
# Price cells (column D) that look like plain numbers are entered with a
# leading apostrophe (classic Excel "force text" prefix) and then the
# cell style is reset to "Normal" so the value stays text (matching the
# original inlineStr cells) without leaving a stray number-format style
# behind. Percentage cells (column E) already contain spaces/"%" so they
# round-trip as text without any extra handling.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.602.78"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "1.658.48"
$ws.Range("E3").Value = "  -4.24%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("E6").Value = "  -2.84%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'24.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E9").Value = "  -2.65%  "
$ws.Range("D11").Value = "'0.0879"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "1.892.98"
$ws.Range("E12").Value = "  -4.25%  "
$ws.Range("D13").Value = "1.701.45"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "'4.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").Value = "'0.567"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").Value = "27.584.93"
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("D18").Value = "'240.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("E19").Value = "  -3.50%  "
$ws.Range("E20").Value = "  -4.46%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'4.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("D23").Value = "'9.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.07%  "
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("D25").Value = "'146.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("E26").Value = "  -4.64%  "
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").Value = "1.458.93"
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("E35").Value = "  -4.91%  "
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").Value = "'0.926"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.49%  "
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").Value = "'0.572"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.17%  "
$ws.Range("D40").Value = "'69.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("E41").Value = "  -4.32%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'5.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "1.801.39"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "'88.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("E49").Value = "  -6.18%  "
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.91%  "
